# "add back the quest config"
# A quest-config row ("传记" / Biography, Id=6, Icon=MainIcon3) that used to
# sit between "成就" (Id=3) and "挑战" (Id=7) was removed at some point; this
# restores it. Re-inserting the row at sheet row 7 pushes every following
# row down by one (Id list, icon names, level/x/y/etc. all travel with their
# row), which reproduces the whole cascade of row-by-row value shifts seen
# in the diff, and the table/autofilter/dimension grow from L20 to L21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 7 ("挑战"), shifting rows 7-20
# down to 8-21.
[void]$ws.Rows.Item(7).Insert()

# Populate the newly blank row 7 with the restored quest-config entry.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "传记"
$ws.Range("C7").Value = "查看自己的传记(T)"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = "MainIcon3"

# Grow the worksheet table ("表1") and its autofilter to cover the new row;
# Excel keeps dimension/table ranges in sync with the data (A1:L20 -> A1:L21).
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A1:L21"))

# Match the author's final cursor position recorded in the sheet view.
[void]$ws.Range("E11").Select()
